$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr2 = New-Object 'object[,]' 1,41
$arr2[0,0] = "Friendly Matches"
$arr2[0,1] = "2025-12-23"
$arr2[0,2] = "09:30:00"
$arr2[0,3] = "Tombense MG"
$arr2[0,4] = "Desportiva"
$arr2[0,5] = 1.32
$arr2[0,6] = 1.81
$arr2[0,7] = 2.44
$arr2[0,8] = 60
$arr2[0,9] = 3.2
$arr2[0,10] = 40
$arr2[0,11] = 1.32
$arr2[0,12] = 1.04
$arr2[0,13] = 1.1
$arr2[0,14] = 1.26
$arr2[0,15] = 1.47
$arr2[0,16] = 1.26
$arr2[0,17] = 1.47
$arr2[0,18] = 1.05
$arr2[0,19] = 1.04
$arr2[0,20] = 1.04
$arr2[0,21] = 1.01
$arr2[0,22] = 2.22
$arr2[0,23] = 1000
$arr2[0,24] = 1000
$arr2[0,25] = 1000
$arr2[0,26] = 1000
$arr2[0,27] = 1000
$arr2[0,28] = 980
$arr2[0,29] = 1000
$arr2[0,30] = 1000
$arr2[0,31] = 1000
$arr2[0,32] = 1000
$arr2[0,33] = 60
$arr2[0,34] = 1000
$arr2[0,35] = 1000
$arr2[0,36] = 1000
$arr2[0,37] = 1000
$arr2[0,38] = 1000
$arr2[0,39] = 29
$arr2[0,40] = 1000
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("A2:AO2").Value = $arr2

# Row 3
$arr3 = New-Object 'object[,]' 1,41
$arr3[0,0] = "Algerian Ligue 1"
$arr3[0,1] = "2025-12-23"
$arr3[0,2] = "15:30:00"
$arr3[0,3] = "MC Alger"
$arr3[0,4] = "ES Ben Aknoun"
$arr3[0,5] = 1.36
$arr3[0,6] = 1.39
$arr3[0,7] = 13
$arr3[0,8] = 16
$arr3[0,9] = 4.7
$arr3[0,10] = 5.3
$arr3[0,11] = 1.5
$arr3[0,12] = 1.09
$arr3[0,13] = 3
$arr3[0,14] = 1.44
$arr3[0,15] = 1.68
$arr3[0,16] = 2.28
$arr3[0,17] = 1.24
$arr3[0,18] = 4.5
$arr3[0,19] = 2.78
$arr3[0,20] = 1.55
$arr3[0,21] = 1.06
$arr3[0,22] = 3.5
$arr3[0,23] = 12
$arr3[0,24] = 30
$arr3[0,25] = 1000
$arr3[0,26] = 1000
$arr3[0,27] = 5.7
$arr3[0,28] = 14
$arr3[0,29] = 65
$arr3[0,30] = 1000
$arr3[0,31] = 6.6
$arr3[0,32] = 11.5
$arr3[0,33] = 130
$arr3[0,34] = 1000
$arr3[0,35] = 11.5
$arr3[0,36] = 20
$arr3[0,37] = 75
$arr3[0,38] = 1000
$arr3[0,39] = 10.5
$arr3[0,40] = 1000
$ws.Range("A3:E3").NumberFormat = "@"
$ws.Range("A3:AO3").Value = $arr3

# Row 4
$arr4 = New-Object 'object[,]' 1,41
$arr4[0,0] = "Bolivian Liga de Futbol Profesional"
$arr4[0,1] = "2025-12-23"
$arr4[0,2] = "16:00:00"
$arr4[0,3] = "Academia de Balompie Boliviano"
$arr4[0,4] = "San Juan FC"
$arr4[0,5] = 1.37
$arr4[0,6] = 1.41
$arr4[0,7] = 3.45
$arr4[0,8] = 7.8
$arr4[0,9] = 5.7
$arr4[0,10] = 1000
$arr4[0,11] = 1.24
$arr4[0,12] = 1.03
$arr4[0,13] = 6
$arr4[0,14] = 1.13
$arr4[0,15] = 2.8
$arr4[0,16] = 1.43
$arr4[0,17] = 1.73
$arr4[0,18] = 2.08
$arr4[0,19] = 1.01
$arr4[0,20] = 1.01
$arr4[0,21] = 1.14
$arr4[0,22] = 3.35
$arr4[0,23] = 1000
$arr4[0,24] = 1000
$arr4[0,25] = 1000
$arr4[0,26] = 1000
$arr4[0,27] = 1000
$arr4[0,28] = 1000
$arr4[0,29] = 1000
$arr4[0,30] = 1000
$arr4[0,31] = 1000
$arr4[0,32] = 1000
$arr4[0,33] = 1000
$arr4[0,34] = 1000
$arr4[0,35] = 1000
$arr4[0,36] = 1000
$arr4[0,37] = 1000
$arr4[0,38] = 1000
$arr4[0,39] = 1000
$arr4[0,40] = 1000
$ws.Range("A4:E4").NumberFormat = "@"
$ws.Range("A4:AO4").Value = $arr4

# Row 5
$arr5 = New-Object 'object[,]' 1,41
$arr5[0,0] = "Friendly Matches"
$arr5[0,1] = "2025-12-23"
$arr5[0,2] = "16:00:00"
$arr5[0,3] = "Serra Branca EC"
$arr5[0,4] = "Maguary"
$arr5[0,5] = 2.12
$arr5[0,6] = 2.5
$arr5[0,7] = 2.88
$arr5[0,8] = 3.9
$arr5[0,9] = 3.25
$arr5[0,10] = 4.2
$arr5[0,11] = 1.37
$arr5[0,12] = 1.07
$arr5[0,13] = 3.6
$arr5[0,14] = 1.27
$arr5[0,15] = 1.9
$arr5[0,16] = 1.81
$arr5[0,17] = 1.4
$arr5[0,18] = 3
$arr5[0,19] = 1.66
$arr5[0,20] = 2.1
$arr5[0,21] = 1.35
$arr5[0,22] = 1.67
$arr5[0,23] = 17.5
$arr5[0,24] = 15.5
$arr5[0,25] = 27
$arr5[0,26] = 70
$arr5[0,27] = 11.5
$arr5[0,28] = 9
$arr5[0,29] = 16
$arr5[0,30] = 44
$arr5[0,31] = 16.5
$arr5[0,32] = 12.5
$arr5[0,33] = 18.5
$arr5[0,34] = 55
$arr5[0,35] = 34
$arr5[0,36] = 26
$arr5[0,37] = 40
$arr5[0,38] = 150
$arr5[0,39] = 19
$arr5[0,40] = 38
$ws.Range("A5:E5").NumberFormat = "@"
$ws.Range("A5:AO5").Value = $arr5

# Row 6
$arr6 = New-Object 'object[,]' 1,41
$arr6[0,0] = "Portuguese Primeira Liga"
$arr6[0,1] = "2025-12-23"
$arr6[0,2] = "17:45:00"
$arr6[0,3] = "Guimaraes"
$arr6[0,4] = "Sporting Lisbon"
$arr6[0,5] = 8.199999999999999
$arr6[0,6] = 8.4
$arr6[0,7] = 1.48
$arr6[0,8] = 1.49
$arr6[0,9] = 4.8
$arr6[0,10] = 5
$arr6[0,11] = 1.39
$arr6[0,12] = 1.06
$arr6[0,13] = 3.95
$arr6[0,14] = 1.33
$arr6[0,15] = 2.04
$arr6[0,16] = 1.94
$arr6[0,17] = 1.39
$arr6[0,18] = 3.4
$arr6[0,19] = 2.14
$arr6[0,20] = 1.86
$arr6[0,21] = 3
$arr6[0,22] = 1.13
$arr6[0,23] = 16.5
$arr6[0,24] = 7.6
$arr6[0,25] = 8.199999999999999
$arr6[0,26] = 12.5
$arr6[0,27] = 24
$arr6[0,28] = 10.5
$arr6[0,29] = 9.4
$arr6[0,30] = 15.5
$arr6[0,31] = 70
$arr6[0,32] = 30
$arr6[0,33] = 29
$arr6[0,34] = 38
$arr6[0,35] = 290
$arr6[0,36] = 140
$arr6[0,37] = 130
$arr6[0,38] = 190
$arr6[0,39] = 180
$arr6[0,40] = 8.199999999999999
$ws.Range("A6:E6").NumberFormat = "@"
$ws.Range("A6:AO6").Value = $arr6

# Row 7
$arr7 = New-Object 'object[,]' 1,41
$arr7[0,0] = "Friendly Matches"
$arr7[0,1] = "2025-12-23"
$arr7[0,2] = "18:00:00"
$arr7[0,3] = "Necaxa"
$arr7[0,4] = "Atletico San Luis"
$arr7[0,5] = 1.93
$arr7[0,6] = 2.28
$arr7[0,7] = 3.45
$arr7[0,8] = 4.5
$arr7[0,9] = 3.3
$arr7[0,10] = 4.3
$arr7[0,11] = 1.35
$arr7[0,12] = 1.07
$arr7[0,13] = 3.4
$arr7[0,14] = 1.32
$arr7[0,15] = 1.92
$arr7[0,16] = 1.7
$arr7[0,17] = 1.36
$arr7[0,18] = 2.78
$arr7[0,19] = 1.67
$arr7[0,20] = 2.08
$arr7[0,21] = 1.29
$arr7[0,22] = 1.78
$arr7[0,23] = 28
$arr7[0,24] = 1000
$arr7[0,25] = 1000
$arr7[0,26] = 1000
$arr7[0,27] = 46
$arr7[0,28] = 19
$arr7[0,29] = 1000
$arr7[0,30] = 1000
$arr7[0,31] = 1000
$arr7[0,32] = 40
$arr7[0,33] = 60
$arr7[0,34] = 1000
$arr7[0,35] = 1000
$arr7[0,36] = 1000
$arr7[0,37] = 1000
$arr7[0,38] = 1000
$arr7[0,39] = 85
$arr7[0,40] = 1000
$ws.Range("A7:E7").NumberFormat = "@"
$ws.Range("A7:AO7").Value = $arr7

# Row 8
$arr8 = New-Object 'object[,]' 1,41
$arr8[0,0] = "Honduras Liga Nacional"
$arr8[0,1] = "2025-12-23"
$arr8[0,2] = "22:00:00"
$arr8[0,3] = "Real Espana"
$arr8[0,4] = "CD Motagua"
$arr8[0,5] = 1.76
$arr8[0,6] = 1.8
$arr8[0,7] = 5.4
$arr8[0,8] = 6.2
$arr8[0,9] = 3.65
$arr8[0,10] = 3.9
$arr8[0,11] = 1.41
$arr8[0,12] = 1.07
$arr8[0,13] = 3.55
$arr8[0,14] = 1.32
$arr8[0,15] = 1.86
$arr8[0,16] = 1.95
$arr8[0,17] = 1.32
$arr8[0,18] = 3.4
$arr8[0,19] = 1.85
$arr8[0,20] = 1.9
$arr8[0,21] = 1.2
$arr8[0,22] = 2.24
$arr8[0,23] = 14
$arr8[0,24] = 19
$arr8[0,25] = 46
$arr8[0,26] = 160
$arr8[0,27] = 8.199999999999999
$arr8[0,28] = 9.4
$arr8[0,29] = 23
$arr8[0,30] = 85
$arr8[0,31] = 10.5
$arr8[0,32] = 10
$arr8[0,33] = 23
$arr8[0,34] = 90
$arr8[0,35] = 19
$arr8[0,36] = 20
$arr8[0,37] = 40
$arr8[0,38] = 1000
$arr8[0,39] = 13
$arr8[0,40] = 110
$ws.Range("A8:E8").NumberFormat = "@"
$ws.Range("A8:AO8").Value = $arr8

# Remove row 9 (Honduras Liga Nacional moved into row 8; original row 9 dropped)
$ws.Rows.Item(9).Delete()

Write-Host "edit complete"